$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows (2-4) down to (3-5)
$ws.Rows.Item(2).Insert()

# The insert operation copies the header row's bold/border formatting down;
# strip that back off so the new row matches the plain style of the other
# data rows.
$ws.Rows.Item(2).ClearFormats()

# Fill in the new row 2 with the latest week's data (same pattern as the
# old row 4 but for the newest date)
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 44630
$ws.Cells.Item(2, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = 100112043
$ws.Cells.Item(2, 7).Value = "Pepino dulce"
$ws.Cells.Item(2, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(2, 11).Value = 15000
$ws.Cells.Item(2, 12).Value = 16000
$ws.Cells.Item(2, 13).Value = 15500
$ws.Cells.Item(2, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(2, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 16).Value = 861
$ws.Cells.Item(2, 17).Value = 18
$ws.Cells.Item(2, 18).Value = "Hortaliza"
